# BurndownChart 20100927 edit
# - Fix the mojibake'd (double UTF-8 encoded) Spanish text in the Sprint
#   sheet and rename the table header row to Title/Weight/Status/Remaining
# - Update the Burndown Chart sheet's first two "Puntos" values (19 -> 18)
# - Switch the active tab from "Sprint" to "Burndown Chart"

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # Sprint
$ws2 = $wb.Worksheets.Item(2)   # Burndown Chart

# --- Sprint sheet: correct text + rename headers -----------------------
# Order matters: it controls the order new shared strings are appended in.
$ws1.Range("A3").Value  = "Crear el diseño general de la master page del sistema SelfManagement"
$ws1.Range("A4").Value  = "Crear el mockup de la pagina de ABM de Campañas para los Jefes de Cuentas"
$ws1.Range("A5").Value  = "Crear el mockup de la pagina de estadisticas globales de las Campañas para los Jefes de Cuentas (utilizando un dashboard y soportando busquedas)"
$ws1.Range("A11").Value = "Diseñar el esquema de la base de datos para el sistema SelfManagement"
$ws1.Range("A12").Value = "Implementar la pantalla de alta de campañas para el sistema SelfManagement"
$ws1.Range("A1").Value  = "Title"
$ws1.Range("B1").Value  = "Weight"
$ws1.Range("C1").Value  = "Status"
$ws1.Range("D1").Value  = "Remaining"

# --- Burndown Chart sheet: update the first two point values -----------
$ws2.Range("B2").Value = 18
$ws2.Range("B3").Value = 18

# --- Selection state -----------------------------------------------------
$ws1.Range("B2:B12").Select() | Out-Null

# --- Make "Burndown Chart" the active tab --------------------------------
$ws2.Activate() | Out-Null
$ws2.Range("A17").Select() | Out-Null
